$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("R1","A6","Y8","A20","K25","P25")
foreach ($c in $cells) {
    $ws.Range($c).Value = "W"
    $ws.Range($c).Interior.Color = 65535
}
